$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "updated marks for 4th week"

# Row 2 (Лазарев / "Исправление Плана управления конфигурациями"): degree of completion 50% -> 90%
$ws.Range("D2").Value = 0.9

# Row 17 (Бурамбекова, mark column): was the shared text "позже" -> now graded as numeric 2
$ws.Range("B17").Value = 2

# Row 19 (Лазарев, mark column): was the shared text "позже" -> now graded "5-"
$ws.Range("B19").Value = "5-"

# Move the active selection to B19, matching the author's final cursor position
$ws.Range("B19").Select()
